$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.636.04"
$ws.Range("E2").Value = "  -4.87%  "
$ws.Range("D3").Value = "1.840.07"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'314.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4225"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.14%  "
$ws.Range("D8").Value = "'0.3624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.05%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'44.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07228"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.75%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.8985"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.18%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'20.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.36%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.782.37"
$ws.Range("E13").Value = "  -7.69%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.557"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.86%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.328"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.61%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.06818"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'77.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.99%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000008902"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.38%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'15.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.18%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.619.77"
$ws.Range("E22").Value = "  -4.90%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.946"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.74%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'10.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.016.03"
$ws.Range("E25").Value = "  -6.71%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.030"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'153.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'18.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.307"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.01%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'110.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.05%  "
$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").Value = "'1.631"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.05%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08884"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.77%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7767"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.80%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.508"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.60%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.070"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.64%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.36%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.082"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.65%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.954"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01917"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.30%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5078"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.58%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1635"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.82%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.725"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.45%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.06639"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'8.220"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.51%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4707"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.37%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'10.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.81%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.623"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.90%  "
